$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value = "[Wen You](https://med.virginia.edu/phs/2019/08/01/you-wen/), [Nichole Szembrot](https://sites.google.com/site/nicholeszembrot/), [Mark Prell](https://www.ers.usda.gov/authors/ers-staff-directory/mark-prell/), [Bruce Weinberg](https://economics.osu.edu/people/weinberg.27)"

$ws.Range("F3").Select()
